$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.801.90'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.639.56'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.81'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.501'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '1.865.87'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = '1.638.35'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.13'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '25.847.01'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.46'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.97'
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.88'
$ws.Range("E24").Value = '  +7.13%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.15'
$ws.Range("E26").Value = '  +2.13%  '
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.93'
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = '1.134.42'
$ws.Range("E37").Value = '  +0.95%  '
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.546'
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.58'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.82'
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("D45").Value = '1.775.09'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.42'
$ws.Range("E50").Value = '  +4.11%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("E51").Value = '  -1.42%  '
